$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81, shifting rows 81:168 down to 82:169
$ws.Rows(81).Insert()

# Populate the new row 81 with the new data point
$ws.Cells.Item(81, 1).Value = 8
$ws.Cells.Item(81, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(81, 3).Value = "Coquimbo"
$ws.Cells.Item(81, 4).Value = 44494
$ws.Cells.Item(81, 5).Value = 4
$ws.Cells.Item(81, 6).Value = 100112012
$ws.Cells.Item(81, 7).Value = "Espinaca"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 2200
$ws.Cells.Item(81, 11).Value = 400
$ws.Cells.Item(81, 12).Value = 500
$ws.Cells.Item(81, 13).Value = 450
$ws.Cells.Item(81, 14).Value = "$/atado 300 a 500 gramos"
$ws.Cells.Item(81, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(81, 16).Value = 900
$ws.Cells.Item(81, 17).Value = 0.5
$ws.Cells.Item(81, 18).Value = "Hortaliza"
